# Updated cryptos list on Thu Aug  1 05:35:35 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to keep a literal text value (matches the source
    # data which stores prices/volumes as inline strings, not numbers),
    # then restore the default "Normal" style so no stray number format
    # is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "64.143.42"
$ws.Range("E2").Value = "  -2.57%  "

Set-TextValue $ws.Range("D3") "3.171.30"
$ws.Range("E3").Value = "  -3.29%  "

$ws.Range("E4").Value = "  +0.02%  "

Set-TextValue $ws.Range("D5") "569.40"
$ws.Range("E5").Value = "  -2.35%  "

Set-TextValue $ws.Range("D6") "169.26"
$ws.Range("E6").Value = "  -5.76%  "

Set-TextValue $ws.Range("D7") "0.607"
$ws.Range("E7").Value = "  -5.38%  "

$ws.Range("E8").Value = "  +0.03%  "

Set-TextValue $ws.Range("D9") "3.172.55"
$ws.Range("E9").Value = "  -3.00%  "

$ws.Range("E10").Value = "  -3.33%  "

$ws.Range("E11").Value = "  -0.77%  "

Set-TextValue $ws.Range("D12") "0.384"
$ws.Range("E12").Value = "  -4.20%  "

Set-TextValue $ws.Range("D13") "3.722.55"
$ws.Range("E13").Value = "  -3.15%  "

$ws.Range("E14").Value = "  -1.76%  "

Set-TextValue $ws.Range("D15") "64.229.64"
$ws.Range("E15").Value = "  -2.56%  "

Set-TextValue $ws.Range("D16") "25.33"
$ws.Range("E16").Value = "  -2.53%  "

$ws.Range("E17").Value = "  -2.30%  "

Set-TextValue $ws.Range("D18") "3.161.11"
$ws.Range("E18").Value = "  -3.57%  "

Set-TextValue $ws.Range("D19") "416.49"
$ws.Range("E19").Value = "  -2.16%  "

Set-TextValue $ws.Range("D20") "12.81"
$ws.Range("E20").Value = "  -2.83%  "

Set-TextValue $ws.Range("D21") "5.34"
$ws.Range("E21").Value = "  -2.66%  "

$ws.Range("E22").Value = "  -3.76%  "

Set-TextValue $ws.Range("D23") "1.00"
$ws.Range("E23").Value = "  -0.12%  "

Set-TextValue $ws.Range("D24") "70.01"
$ws.Range("E24").Value = "  -2.20%  "

Set-TextValue $ws.Range("D25") "0.202"
$ws.Range("E25").Value = "  +2.79%  "

Set-TextValue $ws.Range("D26") "0.491"
$ws.Range("E26").Value = "  -3.50%  "

$ws.Range("E27").Value = "  -5.33%  "

Set-TextValue $ws.Range("D28") "8.76"
$ws.Range("E28").Value = "  -1.22%  "

Set-TextValue $ws.Range("D29") "0.997"
$ws.Range("E29").Value = "  -0.22%  "

$ws.Range("E30").Value = "  -6.06%  "

Set-TextValue $ws.Range("D31") "21.77"
$ws.Range("E31").Value = "  -1.96%  "

$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("E33").Value = "  -2.44%  "

Set-TextValue $ws.Range("D34") "6.34"
$ws.Range("E34").Value = "  -3.53%  "

$ws.Range("E35").Value = "  -4.35%  "

Set-TextValue $ws.Range("D36") "155.23"
$ws.Range("E36").Value = "  -2.30%  "

$ws.Range("E37").Value = "  -3.77%  "

Set-TextValue $ws.Range("D38") "2.698.40"
$ws.Range("E38").Value = "  -3.33%  "

Set-TextValue $ws.Range("D39") "1.70"
$ws.Range("E39").Value = "  -5.36%  "

Set-TextValue $ws.Range("D40") "24.52"
$ws.Range("E40").Value = "  -6.72%  "

$ws.Range("E41").Value = "  -3.24%  "

Set-TextValue $ws.Range("D42") "38.81"
$ws.Range("E42").Value = "  -2.87%  "

Set-TextValue $ws.Range("D43") "0.709"
$ws.Range("E43").Value = "  -7.05%  "

Set-TextValue $ws.Range("D44") "0.0624"
$ws.Range("E44").Value = "  -4.89%  "

Set-TextValue $ws.Range("D45") "5.65"
$ws.Range("E45").Value = "  -3.89%  "

Set-TextValue $ws.Range("D46") "21.81"
$ws.Range("E46").Value = "  -5.28%  "

$ws.Range("E47").Value = "  -2.05%  "

Set-TextValue $ws.Range("D48") "295.37"
$ws.Range("E48").Value = "  -6.14%  "

Set-TextValue $ws.Range("D49") "2.05"
$ws.Range("E49").Value = "  -10.01%  "

Set-TextValue $ws.Range("D50") "1.00"
$ws.Range("E50").Value = "  +0.04%  "

Set-TextValue $ws.Range("D51") "0.0991"
$ws.Range("E51").Value = "  -4.19%  "
